# Generate Report for Handoff
# Adds a new row (for af9db61f-3e0a-4c88-864d-cdcbfd2dc2f7.md) to each of the
# three report sheets: Overview, zh-cn, de-de. Mirrors the existing row for
# 462da0bf-2edc-4240-a168-1c795d3f63f0.md, extends each table by one row, and
# adds a matching hyperlink + relationship on every sheet.

$wb = $excel.ActiveWorkbook

$commitSha = "e2f2665780b0faff9872ebba9d0043a71fb45253"
$newGuid   = "af9db61f-3e0a-4c88-864d-cdcbfd2dc2f7"
$newMd     = "$newGuid.md"
$newPath   = "e2e\$newGuid.md"
$newUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitSha/e2e/$newGuid.md"

$xliffHash = "b29a507db1b5ab11d81d53e7164c137fba20fd16"
$zhXlf     = "$newGuid.$xliffHash.zh-cn.xlf"
$deXlf     = "$newGuid.$xliffHash.de-de.xlf"
$hoDateZh  = "2016-08-22 18:42:41"
$hoDateDe  = "2016-08-22 18:42:46"

# ---------------------------------------------------------------------------
# Sheet "Overview" -> row 3
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = $newMd
$wsOverview.Range("B3").Value = $newPath
$wsOverview.Range("C3").Value = ".md"
# D3 ("Publish URL") mirrors D2, which is blank - leave the cell empty.
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = $hoDateDe
$wsOverview.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $newUrl, "", "", $newPath) | Out-Null

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn" -> row 3
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A3").Value = $newMd
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "False"
$wsZh.Range("G3").Value = $zhXlf
$wsZh.Range("H3").Value = $hoDateZh
$wsZh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
# I3 ("Latest Target File") and J3 ("Latest Handback File") mirror row 2 - blank.
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"
$wsZh.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
# L3 ("Reference Tokens") mirrors row 2 - blank.
$wsZh.Range("M3").Value = "True"
# N3 ("Dependency From") mirrors row 2 - blank.
$wsZh.Range("O3").Value = "False"
# P3 ("Error Detail") mirrors row 2 - blank.

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $newUrl, "", "", $newMd) | Out-Null

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P3"))

# ---------------------------------------------------------------------------
# Sheet "de-de" -> row 3
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A3").Value = $newMd
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "False"
$wsDe.Range("G3").Value = $deXlf
$wsDe.Range("H3").Value = $hoDateDe
$wsDe.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
# I3 ("Latest Target File") and J3 ("Latest Handback File") mirror row 2 - blank.
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDe.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
# L3 ("Reference Tokens") mirrors row 2 - blank.
$wsDe.Range("M3").Value = "True"
# N3 ("Dependency From") mirrors row 2 - blank.
$wsDe.Range("O3").Value = "False"
# P3 ("Error Detail") mirrors row 2 - blank.

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $newUrl, "", "", $newMd) | Out-Null

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P3"))
